$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).Date.AddDays(45212)

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
